# 30% de programa terminado
#
# Adds the "Menor_optima.prn" report sheet (mirroring the Mayor_maxima.prn /
# Pmax_Pgen.prn column layout) and touches up the [MW]/[%] unit suffixes on
# the shared POT_MAX / POT_GEN / MAX_GEN / RESERVA% / PORCENTAJE / RESOPT
# headers used by the existing "*.prn" report sheets.

$wb = $excel.ActiveWorkbook

$wsPmax  = $wb.Worksheets.Item("Pmax_Pgen.prn")
$wsMayor = $wb.Worksheets.Item("Mayor_maxima.prn")

# These two sheets share the same header strings, so relabelling both keeps
# the shared-string table in sync with the rest of the workbook.
foreach ($ws in @($wsPmax, $wsMayor)) {
    $ws.Range("D1").Value = "POT_MAX[MW]"
    $ws.Range("E1").Value = "POT_GEN[MW]"
    $ws.Range("F1").Value = "MAX_GEN[MW]"
    $ws.Range("G1").Value = "RESERVA[%]"
    $ws.Range("H1").Value = "PORCENTAJE[%]"
}

# First pass at renaming the RESOPT column (touches both sheets) ...
foreach ($ws in @($wsPmax, $wsMayor)) {
    $ws.Range("J1").Value = "RES_OPT[[%]"
}

# ... then fix it up, but only on Mayor_maxima.prn.
$wsMayor.Range("J1").Value = "RESOPT[%]"

# Insert the new "Menor_optima.prn" output sheet right before "Reserva.rep",
# reusing the same column headers as Mayor_maxima.prn.
$wsReservaRep = $wb.Worksheets.Item("Reserva.rep")
$wsMenor = $wb.Worksheets.Add($wsReservaRep)
$wsMenor.Name = "Menor_optima.prn"

$wsMenor.Range("A1").Value = "IBUS"
$wsMenor.Range("B1").Value = "NOMBRE"
$wsMenor.Range("C1").Value = "ID"
$wsMenor.Range("D1").Value = "POT_MAX[MW]"
$wsMenor.Range("E1").Value = "POT_GEN[MW]"
$wsMenor.Range("F1").Value = "MAX_GEN[MW]"
$wsMenor.Range("G1").Value = "RESERVA[%]"
$wsMenor.Range("H1").Value = "PORCENTAJE[%]"
$wsMenor.Range("I1").Value = "DATO"
$wsMenor.Range("J1").Value = "RESOPT[%]"
